$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.129.80'
$ws.Range('E2').Value = '  +0.18%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.748.83'
$ws.Range('E3').Value = '  +0.33%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.51'
$ws.Range('E5').Value = '  -0.06%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.98'
$ws.Range('E6').Value = '  -0.55%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.748.01'
$ws.Range('E7').Value = '  +0.34%  '

$ws.Range('E9').Value = '  +1.17%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.170'
$ws.Range('E10').Value = '  +3.21%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.38'
$ws.Range('E11').Value = '  +0.91%  '

$ws.Range('E12').Value = '  +0.11%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.91'
$ws.Range('E13').Value = '  -0.69%  '

$ws.Range('E14').Value = '  +1.42%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.373.84'
$ws.Range('E15').Value = '  +0.31%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.746.80'
$ws.Range('E16').Value = '  +0.27%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.131.38'
$ws.Range('E17').Value = '  +0.28%  '

$ws.Range('E18').Value = '  +1.57%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.40'
$ws.Range('E19').Value = '  +1.15%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.113'
$ws.Range('E20').Value = '  -1.62%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.06'
$ws.Range('E21').Value = '  +8.74%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '492.58'
$ws.Range('E22').Value = '  -1.05%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.727'
$ws.Range('E23').Value = '  +0.67%  '

$ws.Range('E24').Value = '  +7.98%  '

$ws.Range('E25').Value = '  -0.56%  '

$ws.Range('E26').Value = '  -0.21%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.27'
$ws.Range('E27').Value = '  -0.28%  '

$ws.Range('E28').Value = '  -0.83%  '

$ws.Range('E29').Value = '  -0.05%  '

$ws.Range('E30').Value = '  +0.86%  '

$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.13'
$ws.Range('E31').Value = '  +1.84%  '

$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.46'
$ws.Range('E32').Value = '  +1.68%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.893.85'
$ws.Range('E33').Value = '  +0.14%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '31.43'
$ws.Range('E34').Value = '  -0.89%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.681.06'
$ws.Range('E35').Value = '  +0.37%  '

$ws.Range('E36').Value = '  -0.11%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  +0.01%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.01'
$ws.Range('E38').Value = '  +0.19%  '

$ws.Range('E39').Value = '  +2.49%  '

$ws.Range('E40').Value = '  +3.35%  '

$ws.Range('E41').Value = '  +0.20%  '

$ws.Range('E42').Value = '  +5.82%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '48.69'
$ws.Range('E43').Value = '  -0.57%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '424.98'
$ws.Range('E44').Value = '  -2.23%  '

$ws.Range('E45').Value = '  -0.68%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.46'
$ws.Range('E46').Value = '  +0.59%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.09'
$ws.Range('E48').Value = '  -1.03%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.61'
$ws.Range('E49').Value = '  -0.24%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.788.22'
$ws.Range('E50').Value = '  +1.57%  '

$ws.Range('E51').Value = '  +0.03%  '
